$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.105.58"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.655.20"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.73"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5236"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2608"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06351"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.37"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07794"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.503"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.660.95"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5476"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "0.0₅8213"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.39"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "26.125.42"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.582"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.56"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.035"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.17"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1242"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.251"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.13"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05905"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.527"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.251"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.592"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9545"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.787"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5703"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01618"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.821"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8501"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "1.033.73"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.17"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").Value = "1.800.25"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.22"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4302"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.476"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05167"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.845"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09712"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.02%  "
